$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1711
$ws.Range("F4").Value = 1123
$ws.Range("F6").Value = 141
$ws.Range("F7").Value = 1398
$ws.Range("F9").Value = 91
$ws.Range("F10").Value = 615
$ws.Range("F11").Value = 128
$ws.Range("F12").Value = 84
$ws.Range("F14").Value = 453
$ws.Range("F15").Value = 463
$ws.Range("F16").Value = 126
$ws.Range("C17").Value = '上海·元宵AuPoRo音乐动漫FES'
$ws.Range("F18").Value = 680
$ws.Range("F19").Value = 2537
$ws.Range("F21").Value = 45
$ws.Range("F22").Value = 15
$ws.Range("F24").Value = 281
$ws.Range("F25").Value = 170
$ws.Range("F26").Value = 8
$ws.Range("F28").Value = 560
$ws.Range("F29").Value = 915
$ws.Range("F31").Value = 55
$ws.Range("F34").Value = 26
$ws.Range("F35").Value = 238

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 715
$ws.Range("F5").Value = 605
$ws.Range("G5").Value = 280
$ws.Range("F6").Value = 605
$ws.Range("G6").Value = 280
$ws.Range("F8").Value = 1
$ws.Range("F12").Value = 268
$ws.Range("F15").Value = 322
$ws.Range("F16").Value = 322
$ws.Range("F17").Value = 66
$ws.Range("F19").Value = 926
$ws.Range("F24").Value = 24
$ws.Range("F26").Value = 225
$ws.Range("F27").Value = 219
$ws.Range("F29").Value = 185

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("G5").Value = '不可售'
$ws.Range("F6").Value = 2259
$ws.Range("F7").Value = 895
$ws.Range("F10").Value = 1103
$ws.Range("F11").Value = 236
$ws.Range("F12").Value = 69

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2259
$ws.Range("F5").Value = 1711
$ws.Range("F8").Value = 895
$ws.Range("F9").Value = 1103
$ws.Range("F10").Value = 236
$ws.Range("F11").Value = 69
$ws.Range("F12").Value = 715
$ws.Range("F13").Value = 1123
$ws.Range("F15").Value = 141
$ws.Range("F16").Value = 1398
$ws.Range("F17").Value = 605
$ws.Range("G17").Value = 280
$ws.Range("F19").Value = 91
$ws.Range("F20").Value = 615
$ws.Range("F21").Value = 128
$ws.Range("F23").Value = 84
$ws.Range("F25").Value = 453
$ws.Range("F26").Value = 463
$ws.Range("C27").Value = '上海·元宵AuPoRo音乐动漫FES'
$ws.Range("F28").Value = 680
$ws.Range("F29").Value = 2537
$ws.Range("F31").Value = 15
$ws.Range("F32").Value = 281
$ws.Range("F33").Value = 268
$ws.Range("F34").Value = 170
$ws.Range("F37").Value = 560
$ws.Range("F38").Value = 915
$ws.Range("F39").Value = 322
$ws.Range("F40").Value = 66
$ws.Range("F42").Value = 55
$ws.Range("F45").Value = 24
$ws.Range("F46").Value = 225
$ws.Range("F47").Value = 219
$ws.Range("F48").Value = 185
$ws.Range("F50").Value = 26
$ws.Range("F51").Value = 238
